# Moved all Documents to docs folder and added ERD and UML Diagram
$wb = $excel.ActiveWorkbook

# New date value (2024-03-28 -> serial 45379) and text for the new row
$newDate   = (Get-Date -Year 2024 -Month 3 -Day 28).Date
$newHours  = "2h 20min"
$newDesc   = "Klassendiagramm und ERD erstellt"

$sheetNames = @("PaulSchein", "ReneMifka", "LukasPerger")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)

    # find first empty row in column A (after the header + existing data rows)
    $lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row  # xlUp = -4162
    $newRow = $lastRow + 1

    # Copy the cell above first so the new cell inherits the same style (s="1",
    # numFmtId 14 date format) instead of Excel fabricating a brand-new style.
    $dateCell = $ws.Cells.Item($newRow, 1)
    $ws.Cells.Item($newRow - 1, 1).Copy($dateCell) | Out-Null
    $dateCell.Value = $newDate

    $ws.Cells.Item($newRow, 2).Value = $newHours
    $ws.Cells.Item($newRow, 3).Value = $newDesc
}

# --- Update each sheet's view (selection / zoom / active tab) ---

# PaulSchein: selection moves down onto the newly added row A4:C4
$wsPaul = $wb.Worksheets.Item("PaulSchein")
$wsPaul.Activate()
$wsPaul.Range("A4:C4").Select()

# ReneMifka: becomes the active tab, zoom bumped to 325%, selection at C9
$wsRene = $wb.Worksheets.Item("ReneMifka")
$wsRene.Activate()
$wsRene.Application.ActiveWindow.Zoom = 325
$wsRene.Range("C9").Select()

# LukasPerger: no longer the active tab, selection moves to C8
$wsLukas = $wb.Worksheets.Item("LukasPerger")
$wsLukas.Activate()
$wsLukas.Range("C8").Select()

# ReneMifka (workbook activeTab index 1) is the tab left active/selected
$wsRene.Activate()
